$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixText = @'
# Fix NERDTree not working on cygwin:
* __Symptom__ when open file works but not able to open directory, And the expand/collapse arrow cannot show 
* __Reason__ the failure of showing the arrow cause the direcotry select error
* __Solution__ override the arrow charactor in vimrc by adding below lines:
```
let g:NERDTreeDirArrowExpandable="+"
let g:NERDTreeDirArrowCollapsible="-"
```
'@

$ws.Range("A25").Value = "nerdtree"
$ws.Range("B25").Value = "debug"
$ws.Range("C25").Value = $fixText

$ws.Range("C25").WrapText = $true
$ws.Rows(25).RowHeight = 135

# Mirror Excel's post-entry UI state: selection drops to the next row,
# and the view scrolls so the new row is visible near the bottom.
$ws.Range("C26").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1
